$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "71.281.11"
$ws.Range("E2").Value2 = "  +4.97%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.621.81"

# Row 4
$ws.Range("E4").Value2 = "  +0.00%  "

# Row 5
Set-TextValue $ws.Range("D5") "605.55"
$ws.Range("E5").Value2 = "  +2.92%  "

# Row 6
Set-TextValue $ws.Range("D6") "180.23"
$ws.Range("E6").Value2 = "  +3.41%  "

# Row 7
$ws.Range("E7").Value2 = "  -0.05%  "

# Row 8
$ws.Range("E8").Value2 = "  +1.76%  "

# Row 9
Set-TextValue $ws.Range("D9") "2.621.18"
$ws.Range("E9").Value2 = "  +5.34%  "

# Row 10
$ws.Range("E10").Value2 = "  +15.10%  "

# Row 11
$ws.Range("E11").Value2 = "  +0.52%  "

# Row 12
$ws.Range("E12").Value2 = "  +4.01%  "

# Row 13
Set-TextValue $ws.Range("D13") "5.03"
$ws.Range("E13").Value2 = "  +1.37%  "

# Row 15
Set-TextValue $ws.Range("D15") "26.63"
$ws.Range("E15").Value2 = "  +4.98%  "

# Row 16
$ws.Range("E16").Value2 = "  +8.01%  "

# Row 17
Set-TextValue $ws.Range("D17") "71.152.44"
$ws.Range("E17").Value2 = "  +5.01%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.633.55"
$ws.Range("E18").Value2 = "  +6.95%  "

# Row 19
Set-TextValue $ws.Range("D19") "379.06"
$ws.Range("E19").Value2 = "  +9.32%  "

# Row 20
$ws.Range("E20").Value2 = "  +7.02%  "

# Row 21
Set-TextValue $ws.Range("D21") "11.48"
$ws.Range("E21").Value2 = "  +5.88%  "

# Row 22
$ws.Range("E22").Value2 = "  +0.80%  "

# Row 23
Set-TextValue $ws.Range("D23") "71.91"
$ws.Range("E23").Value2 = "  +1.56%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.43"
$ws.Range("E24").Value2 = "  +6.40%  "

# Row 25
$ws.Range("E25").Value2 = "  +0.08%  "

# Row 26
Set-TextValue $ws.Range("D26") "1.84"
$ws.Range("E26").Value2 = "  +8.10%  "

# Row 27
$ws.Range("E27").Value2 = "  +7.93%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.755.08"
$ws.Range("E28").Value2 = "  +5.38%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value2 = "  +0.34%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0₃0953"
$ws.Range("E30").Value2 = "  +7.20%  "

# Row 31
Set-TextValue $ws.Range("D31") "531.42"
$ws.Range("E31").Value2 = "  +6.75%  "

# Row 32
$ws.Range("E32").Value2 = "  +3.30%  "

# Row 33
$ws.Range("E33").Value2 = "  +5.58%  "

# Row 34
$ws.Range("E34").Value2 = "  +3.76%  "

# Row 35
$ws.Range("E35").Value2 = "  +0.03%  "

# Row 36
Set-TextValue $ws.Range("D36") "165.50"
$ws.Range("E36").Value2 = "  +0.76%  "

# Row 37
$ws.Range("E37").Value2 = "  -0.91%  "

# Row 38
Set-TextValue $ws.Range("D38") "19.13"
$ws.Range("E38").Value2 = "  +5.07%  "

# Row 39
$ws.Range("B39").Value2 = "Stacks"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D39") "1.88"
$ws.Range("E39").Value2 = "  +8.76%  "

# Row 40
$ws.Range("B40").Value2 = "WhiteBITCoin"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D40") "18.98"
$ws.Range("E40").Value2 = "  +1.88%  "

# Row 41
$ws.Range("E41").Value2 = "  +4.87%  "

# Row 42
$ws.Range("E42").Value2 = "  +0.02%  "

# Row 43
$ws.Range("B43").Value2 = "dogwifhat"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.60"
$ws.Range("E43").Value2 = "  +9.72%  "

# Row 44
$ws.Range("B44").Value2 = "RenderToken"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws.Range("D44") "5.02"
$ws.Range("E44").Value2 = "  +5.18%  "

# Row 45
$ws.Range("E45").Value2 = "  +2.83%  "

# Row 47
Set-TextValue $ws.Range("D47") "153.46"
$ws.Range("E47").Value2 = "  +3.35%  "

# Row 48
Set-TextValue $ws.Range("D48") "3.66"
$ws.Range("E48").Value2 = "  +3.66%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0₆0266"
$ws.Range("E49").Value2 = "  +5.71%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.532"
$ws.Range("E50").Value2 = "  +3.87%  "

# Row 51
$ws.Range("E51").Value2 = "  +7.51%  "
